$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49. This shifts rows 49..113 down to 50..114
# and automatically updates the sheet dimension to A1:R114.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly data point.
$ws.Range("A49").Value = 6
$ws.Range("B49").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C49").Value = "Metropolitana"
$ws.Range("D49").Value = 44477
$ws.Range("D49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E49").Value = 13
$ws.Range("F49").Value = 100112029
$ws.Range("G49").Value = "Orégano"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 35
$ws.Range("K49").Value = 8500
$ws.Range("L49").Value = 9000
$ws.Range("M49").Value = 8729
$ws.Range("N49").Value = "`$/docena de atados"
$ws.Range("O49").Value = "Región Metropolitana"
$ws.Range("P49").Value = 2910
$ws.Range("Q49").Value = 3
$ws.Range("R49").Value = "Hortaliza"
